# Adding Header field in to 1PProfile
#
# The STATUS column (L) no longer carries a hard-coded "PASS" value for
# the existing test rows (2-6) -- clear it out.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("L2:L6").ClearContents() | Out-Null

# Leave the grid scrolled/selected the way the author left it: focus on
# the (now-empty) STATUS column.
$ws.Range("L2:L6").Select() | Out-Null
$excel.ActiveWindow.ScrollColumn = 5
$excel.ActiveWindow.ScrollRow = 1
